# Apply the "Finalized Experiments with Participant Generation" edit:
# - rename the 5 task-order sheets with refreshed timestamp ids
# - update generated stim/trial filenames (and swap the eyes open/closed
#   order on the RS sheet) to match the newly generated run

$wb = $excel.ActiveWorkbook

# --- rename sheets -------------------------------------------------------
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912841939447"
$wb.Worksheets.Item(2).Name = "NB_TO-1650291287560992"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912875629926"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912876377008"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291287736451"

# --- sheet 1 (GNG) column B values ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912841621745.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912841769805.csv"
$ws1.Range("B4").Value = "go_stims-16502912841779866.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912841929436.csv"

# --- sheet 2 (NB) column B values -----------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16502912854369676.csv"
$ws2.Range("B3").Value = "ZB-match_7-1650291284602591.csv"
$ws2.Range("B4").Value = "OB-1650291285789912.csv"
$ws2.Range("B5").Value = "TB-16502912871780229.csv"
$ws2.Range("B6").Value = "TB-16502912862272243.csv"
$ws2.Range("B7").Value = "ZB-match_7-1650291284498439.csv"
$ws2.Range("B8").Value = "ZB-match_9-16502912847444227.csv"
$ws2.Range("B9").Value = "OB-1650291285385154.csv"
$ws2.Range("B10").Value = "TB-1650291287542437.csv"

# --- sheet 3 (RS) column B values (order swapped) -------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- sheet 4 (TOL) column B values -----------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912875931566.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912875666442.csv"
$ws4.Range("B4").Value = "MM_stims-1650291287623221.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291287594151.csv"
$ws4.Range("B6").Value = "MM_stims-16502912876367004.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912876242204.csv"

# --- sheet 5 (vSAT) column B values ----------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502912876714535.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912876429887.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502912876931376.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650291287724178.csv"
